$d = $word.ActiveDocument
foreach ($p in $d.Paragraphs) {
    $p.Format.ContextualSpacing = $false
}
